# Merge Jacob's review comments about Magnus into the assessment sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Peer  and self assessment")

# Criterion 1 Online collaboration section - Jacob's row (row 3)
$ws.Range("B3").Value = "Good"
$ws.Range("C3").Value = "1) If asked to communicate, Magnus responds quickly.`n2) Good communicator in group meetings. `n3) Has not been very active in the technical discussions or documenting much`nof his research."
$ws.Range("C3").WrapText = $true

# Criterion 1 International Collaboration section - Jacob's row (row 16)
$ws.Range("B16").Value = "Sufficient"
$ws.Range("C16").Value = "1) Has been participating in meetings and collaboration`n2) Has not been very active with comments/critique. `n"
$ws.Range("C16").WrapText = $true

# Update view: scroll down so row 25 is at the top, move selection to C6
$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Range("C6").Select()
